$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("FEVRIER ")
Write-Output "Found: $found"
Write-Output "Start: $($r.Start) End: $($r.End)"
$s = $r.Start
$e = $r.End
$d.Bookmarks.Add("_GoBack", $d.Range($s, $s))
$r.Delete()
Write-Output "Bookmark added, deletion done"
